$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original data occupied rows 2-21 (20 rows x 3 cols). We need to shift
# that block down by 8 rows (to rows 10-29) to make room for 8 new rows
# inserted right after the header, and then also append 2 brand-new rows
# at the very end (rows 30-31).
#
# We avoid Range.Insert()/EntireRow.Insert() here because it copies the
# formatting of the row above down into the newly created rows (creating
# a spurious extra cell style), which the source edit did not do. Instead
# we just move the existing values down directly and then fill in both
# the new leading rows and the new trailing rows - no shifting machinery,
# no incidental style churn.

$shift = 8
$firstDataRow = 2
$lastDataRow = 21

# Copy from the bottom up so we never overwrite a source row before it has
# been read.
for ($r = $lastDataRow; $r -ge $firstDataRow; $r--) {
    $destRow = $r + $shift
    $ws.Cells.Item($destRow, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($destRow, 2).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($destRow, 3).Value = $ws.Cells.Item($r, 3).Value2
}

# New rows inserted right after the header row (now rows 2-9).
$newTopRows = @(
    @(-0.0042760567739605, -0.1111774742603302, -0.0980438739061355),
    @(0.0734565481543541, 0.1905899494886398, -0.1611157059669494),
    @(-0.101709060370922, -0.09025534242391579, 0.1893682330846786),
    @(-0.0087048299610614, 0.024892758578062, 0.00167987938039),
    @(0.0125227374956011, -0.0652098655700683, 0.066737025976181),
    @(-0.0397062413394451, 0.0247400421649217, 0.030695978552103),
    @(0.012980886735022, 0.0633772686123848, -0.0366519130766391),
    @(-0.0360410511493682, 0.0274889357388019, -0.0474947728216648)
)

for ($i = 0; $i -lt $newTopRows.Count; $i++) {
    $r = $firstDataRow + $i
    $vals = $newTopRows[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
}

# Two additional rows appended at the bottom of the table (new rows 30-31).
$newBottomRows = @(
    @(-0.0519235469400882, -0.107512280344963, 0.1065959855914115),
    @(0.0108428578823804, -0.0210748501121997, -0.0189368221908807)
)

$lastRow = $lastDataRow + $shift
for ($i = 0; $i -lt $newBottomRows.Count; $i++) {
    $r = $lastRow + 1 + $i
    $vals = $newBottomRows[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
}
